$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 2 (shifting SELECT rows down to 4:5)
$ws.Rows("2:3").Insert(-4121)

# Clear inherited header formatting on the newly inserted rows
$ws.Range("A2:G3").Style = "Normal"

# New row 2: CHOICE / SYMBOL_2000_CHOICE variant
$ws.Range("A2").Value = "AAA_CSA"
$ws.Range("B2").Value = "CHOICE"
$ws.Range("C2").Value = "CA"
$ws.Range("D2").Value = "SYMBOL_2000_CHOICE"
$ws.Range("E2").Value = 20000101
$ws.Range("F2").Value = 20180729
$ws.Range("G2").Value = "MSRP_2000_CHOICE"

# New row 3: CHOICE / SYMBOL_2018_CHOICE variant
$ws.Range("A3").Value = "AAA_CSA"
$ws.Range("B3").Value = "CHOICE"
$ws.Range("C3").Value = "CA"
$ws.Range("D3").Value = "SYMBOL_2018_CHOICE"
$ws.Range("E3").Value = 20180730
$ws.Range("F3").Value = 99999999
$ws.Range("G3").Value = "MSRP_2000_CHOICE"

# Existing SELECT rows (now rows 4 and 5) get updated date values
$ws.Range("F4").Value = 20180729

$ws.Range("D5").Value = "SYMBOL_2018"
$ws.Range("E5").Value = 20180730
$ws.Range("F5").Value = 99999999

# Apply the bordered/highlighted date-column formatting (as used by F4)
# to the new cells and the previously-unformatted date cells
$ws.Range("F4").Copy()
$ws.Range("E2:F3").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E5:F5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F7").Select()
